# Updated cryptos list (GitHub Actions data refresh): refresh Price (D) and
# Volume(1h) (E) columns for each coin row; FirstDigitalUSD/ImmutableX swap
# ranking positions (rows 29 <-> 30). D-column values are prefixed with a
# leading apostrophe so Excel keeps them as literal text (preserving
# trailing zeros / thousands-dot separators) instead of auto-coercing them
# into numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'63.381.48"
$ws.Range('E2').Value = '  -4.19%  '
$ws.Range('D3').Value = "'3.094.35"
$ws.Range('E3').Value = '  -4.67%  '
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').Value = "'547.97"
$ws.Range('E5').Value = '  -4.78%  '
$ws.Range('D6').Value = "'136.84"
$ws.Range('E6').Value = '  -11.24%  '
$ws.Range('D7').Value = "'1.00"
$ws.Range('E7').Value = '  +0.12%  '
$ws.Range('D8').Value = "'3.085.26"
$ws.Range('E8').Value = '  -4.75%  '
$ws.Range('D9').Value = "'0.497"
$ws.Range('E9').Value = '  -2.97%  '
$ws.Range('D10').Value = "'0.156"
$ws.Range('E10').Value = '  -5.43%  '
$ws.Range('D11').Value = "'6.25"
$ws.Range('E11').Value = '  -11.90%  '
$ws.Range('E12').Value = '  -4.13%  '
$ws.Range('D13').Value = "'35.51"
$ws.Range('E13').Value = '  -6.04%  '
$ws.Range('D14').Value = "'0.0000217"
$ws.Range('E14').Value = '  -7.99%  '
$ws.Range('D15').Value = "'3.593.89"
$ws.Range('E15').Value = '  -4.31%  '
$ws.Range('D16').Value = "'63.405.94"
$ws.Range('E16').Value = '  -4.31%  '
$ws.Range('E17').Value = '  -3.05%  '
$ws.Range('D18').Value = "'3.099.81"
$ws.Range('E18').Value = '  -4.50%  '
$ws.Range('D19').Value = "'6.76"
$ws.Range('E19').Value = '  -5.13%  '
$ws.Range('D20').Value = "'489.34"
$ws.Range('E20').Value = '  -12.45%  '
$ws.Range('D21').Value = "'13.65"
$ws.Range('E21').Value = '  -5.53%  '
$ws.Range('D22').Value = "'0.719"
$ws.Range('E22').Value = '  -3.45%  '
$ws.Range('D23').Value = "'7.24"
$ws.Range('E23').Value = '  -8.26%  '
$ws.Range('D24').Value = "'79.17"
$ws.Range('E24').Value = '  -3.59%  '
$ws.Range('D25').Value = "'12.39"
$ws.Range('E25').Value = '  -8.82%  '
$ws.Range('D26').Value = "'0.999"
$ws.Range('E26').Value = '  -0.01%  '
$ws.Range('D27').Value = "'8.46"
$ws.Range('E27').Value = '  -9.78%  '
$ws.Range('D28').Value = "'2.75"
$ws.Range('E28').Value = '  -7.16%  '
$ws.Range('B29').Value = 'FirstDigitalUSD'
$ws.Range('C29').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D29').Value = "'1.00"
$ws.Range('E29').Value = '  -0.14%  '
$ws.Range('B30').Value = 'ImmutableX'
$ws.Range('C30').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D30').Value = "'1.98"
$ws.Range('E30').Value = '  -12.16%  '
$ws.Range('D31').Value = "'26.61"
$ws.Range('E31').Value = '  -4.50%  '
$ws.Range('D32').Value = "'1.13"
$ws.Range('E32').Value = '  -4.13%  '
$ws.Range('D33').Value = "'2.51"
$ws.Range('E33').Value = '  -9.67%  '
$ws.Range('D34').Value = "'58.10"
$ws.Range('E34').Value = '  +4.63%  '
$ws.Range('D35').Value = "'514.17"
$ws.Range('E35').Value = '  -10.23%  '
$ws.Range('D36').Value = "'6.01"
$ws.Range('E36').Value = '  -6.10%  '
$ws.Range('D37').Value = "'5.13"
$ws.Range('E37').Value = '  -11.16%  '
$ws.Range('D38').Value = "'0.0403"
$ws.Range('E38').Value = '  -11.67%  '
$ws.Range('D39').Value = "'3.156.53"
$ws.Range('E39').Value = '  +0.08%  '
$ws.Range('D40').Value = "'0.0803"
$ws.Range('E40').Value = '  -7.28%  '
$ws.Range('D41').Value = "'0.119"
$ws.Range('E41').Value = '  -7.05%  '
$ws.Range('D42').Value = "'8.16"
$ws.Range('E42').Value = '  -5.37%  '
$ws.Range('D43').Value = "'2.66"
$ws.Range('E43').Value = '  -12.91%  '
$ws.Range('E44').Value = '  -5.50%  '
$ws.Range('E45').Value = '  +0.04%  '
$ws.Range('D46').Value = "'2.06"
$ws.Range('E46').Value = '  -10.64%  '
$ws.Range('D47').Value = "'25.24"
$ws.Range('E47').Value = '  -5.17%  '
$ws.Range('D48').Value = "'121.03"
$ws.Range('E48').Value = '  -3.18%  '
$ws.Range('E49').Value = '  -3.97%  '
$ws.Range('D50').Value = "'0.0₃0503"
$ws.Range('E50').Value = '  -9.56%  '
$ws.Range('D51').Value = "'2.03"
$ws.Range('E51').Value = '  -9.94%  '
